$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename header cells: "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404"
# ---------------------------------------------------------------------------
$used = $ws.UsedRange
$lastCol = $used.Columns.Count
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Text
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2310"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2404"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Freeze the header row (split below row 1, scrolling pane starts at A2)
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3) Turn the header+data range into an Excel Table ("ListObject") so the
#    renamed headers above also drive the table's column metadata.
#    Stash the header row's current formatting first and restore it verbatim
#    after the table is created, since Add() otherwise captures whatever
#    formatting the header already has as a one-off conditional-format
#    override on the table definition.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("W1:AQ1")
$headerRange.Copy()
$scratch.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$headerRange.Style = "Normal"

$tableRange = $ws.Range("A1:U64")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$scratch.Clear()
